$d = $word.ActiveDocument

function Replace-WithXml($searchText, $innerXml) {
    $docLenBefore = $d.Content.End
    $find = $d.Content
    $find.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $oldStart = $find.Start
    $oldEnd = $find.End

    $insertPoint = $d.Range($oldStart, $oldStart)
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $innerXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $insertPoint.InsertXML($xml)

    $docLenAfter = $d.Content.End
    $delta = $docLenAfter - $docLenBefore
    $oldRangeNow = $d.Range($oldStart + $delta, $oldEnd + $delta)
    $oldRangeNow.Delete()
}

# --- Change 1: split the "<Proposal Description>" paragraph into two,
#     and move the _GoBack bookmark onto the new paragraph ---
$find1 = $d.Content
$find1.Find.Execute("<Proposal Description>", $false, $false, $false, $false, $false, $true, 1, $false, "^p<Proposal Description>", 2)

$find1b = $d.Content
$find1b.Find.Execute("<Proposal Description>", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$find1b.Font.Italic = $true
$find1b.Font.ItalicBi = $true

$d.Bookmarks.ShowHidden = $true
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$find1c = $d.Content
$find1c.Find.Execute("<Proposal Description>", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$find1c.Collapse(1)
$d.Bookmarks.Add("_GoBack", $find1c)

# --- Change 2: split the compliance sentence into three runs around a
#     grammar-checker proofing mark ("is in compliance with") ---
$xml2 = '<w:r w:rsidRPr="00357271"><w:rPr><w:rFonts w:ascii="Source Sans Pro Light" w:hAnsi="Source Sans Pro Light"/></w:rPr><w:t xml:space="preserve">Having considered the submitted document I confirm it accords with relevant standards and guidance and that it </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Source Sans Pro Light" w:hAnsi="Source Sans Pro Light"/></w:rPr><w:t>is in compliance with</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Source Sans Pro Light" w:hAnsi="Source Sans Pro Light"/></w:rPr><w:t xml:space="preserve"> the condition.</w:t></w:r>'
Replace-WithXml "Having considered the submitted document I confirm it accords with relevant standards and guidance and that it is in compliance with the condition." $xml2

# --- Change 3: drop the stale lastRenderedPageBreak marker ---
$xml3 = '<w:r w:rsidRPr="00357271"><w:rPr><w:rFonts w:ascii="Source Sans Pro Light" w:hAnsi="Source Sans Pro Light"/></w:rPr><w:t>This response relates solely to archaeological issues.</w:t></w:r>'
Replace-WithXml "This response relates solely to archaeological issues." $xml3
